# Weekly fruit/vegetable data refresh: a new week's price observation is
# inserted as a new row at position 196 (pushing every subsequent row down
# by one), growing the used range from A1:R314 to A1:R315.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the new observation.
$ws.Rows.Item(196).EntireRow.Insert()

# Populate the newly inserted row with the new weekly data point.
$ws.Range("A196").Value = 9
$ws.Range("B196").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C196").Value = "Metropolitana"
$ws.Range("D196").Value = 44806
$ws.Range("E196").Value = 13
$ws.Range("F196").Value = 300000001
$ws.Range("G196").Value = "Rabanito"
$ws.Range("H196").Value = "Sin especificar"
$ws.Range("I196").Value = "Primera"
$ws.Range("J196").Value = 7900
$ws.Range("K196").Value = 2500
$ws.Range("L196").Value = 3000
$ws.Range("M196").Value = 2750
$ws.Range("N196").Value = "$/cien unidades (volumen en unidades)"
$ws.Range("O196").Value = "Provincia de Chacabuco"
$ws.Range("P196").Value = 28
$ws.Range("Q196").Value = 100
$ws.Range("R196").Value = "Hortaliza"
